# Update the "custo do toner por periodo de tempo" cost table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B updates
$ws.Range("B2").Value  = 9
$ws.Range("B4").Value  = 29
$ws.Range("B8").Value  = 5
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 166
$ws.Range("B12").Value = 22
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 366
$ws.Range("B15").Value = 14

# Column C updates
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 6
$ws.Range("C12").Value = 2
$ws.Range("C13").Value = 2
$ws.Range("C14").Value = 6
